$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styles for cells whose new value looks numeric,
# so we can force them to remain text (matching the source inline strings)
# without altering their visual style once the true value is restored.
$origStyles = @{}
$origStyles["5_4"] = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$origStyles["6_4"] = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$origStyles["9_4"] = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$origStyles["10_4"] = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$origStyles["12_4"] = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = "@"
$origStyles["14_4"] = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$origStyles["18_4"] = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = "@"
$origStyles["20_4"] = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$origStyles["21_4"] = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = "@"
$origStyles["23_4"] = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$origStyles["24_4"] = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$origStyles["26_4"] = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$origStyles["31_4"] = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$origStyles["32_4"] = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$origStyles["33_4"] = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$origStyles["34_4"] = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$origStyles["35_4"] = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$origStyles["36_4"] = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = "@"
$origStyles["38_4"] = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = "@"
$origStyles["39_4"] = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = "@"
$origStyles["40_4"] = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = "@"
$origStyles["43_4"] = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = "@"
$origStyles["44_4"] = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$origStyles["45_4"] = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$origStyles["46_4"] = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$origStyles["47_4"] = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$origStyles["48_4"] = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$origStyles["51_4"] = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"

# Apply all cell updates
$ws.Cells.Item(2, 4).Value2 = '51.583.35'
$ws.Cells.Item(2, 5).Value2 = '  +1.08%  '
$ws.Cells.Item(3, 4).Value2 = '3.017.19'
$ws.Cells.Item(3, 5).Value2 = '  +2.04%  '
$ws.Cells.Item(4, 5).Value2 = '  +0.05%  '
$ws.Cells.Item(5, 4).Value2 = '378.43'
$ws.Cells.Item(5, 5).Value2 = '  -0.43%  '
$ws.Cells.Item(6, 4).Value2 = '102.30'
$ws.Cells.Item(6, 5).Value2 = '  -0.18%  '
$ws.Cells.Item(7, 5).Value2 = '  +0.29%  '
$ws.Cells.Item(8, 5).Value2 = '  -0.01%  '
$ws.Cells.Item(9, 4).Value2 = '0.589'
$ws.Cells.Item(9, 5).Value2 = '  +0.44%  '
$ws.Cells.Item(10, 4).Value2 = '36.67'
$ws.Cells.Item(10, 5).Value2 = '  +0.66%  '
$ws.Cells.Item(11, 5).Value2 = '  -0.18%  '
$ws.Cells.Item(12, 4).Value2 = '0.0862'
$ws.Cells.Item(12, 5).Value2 = '  +1.24%  '
$ws.Cells.Item(13, 4).Value2 = '3.494.63'
$ws.Cells.Item(13, 5).Value2 = '  +1.95%  '
$ws.Cells.Item(14, 4).Value2 = '18.39'
$ws.Cells.Item(14, 5).Value2 = '  -0.16%  '
$ws.Cells.Item(15, 5).Value2 = '  -0.50%  '
$ws.Cells.Item(16, 4).Value2 = '3.019.88'
$ws.Cells.Item(16, 5).Value2 = '  +2.14%  '
$ws.Cells.Item(17, 5).Value2 = '  -4.08%  '
$ws.Cells.Item(18, 4).Value2 = '10.62'
$ws.Cells.Item(18, 5).Value2 = '  -14.91%  '
$ws.Cells.Item(19, 4).Value2 = '51.547.43'
$ws.Cells.Item(19, 5).Value2 = '  +0.89%  '
$ws.Cells.Item(20, 4).Value2 = '3.10'
$ws.Cells.Item(20, 5).Value2 = '  +0.80%  '
$ws.Cells.Item(21, 4).Value2 = '12.42'
$ws.Cells.Item(21, 5).Value2 = '  +0.31%  '
$ws.Cells.Item(22, 5).Value2 = '  +0.73%  '
$ws.Cells.Item(23, 4).Value2 = '69.89'
$ws.Cells.Item(23, 5).Value2 = '  +0.32%  '
$ws.Cells.Item(24, 4).Value2 = '266.88'
$ws.Cells.Item(24, 5).Value2 = '  -0.24%  '
$ws.Cells.Item(25, 5).Value2 = '  -7.18%  '
$ws.Cells.Item(26, 4).Value2 = '8.28'
$ws.Cells.Item(26, 5).Value2 = '  +3.71%  '
$ws.Cells.Item(27, 5).Value2 = '  +8.69%  '
$ws.Cells.Item(28, 5).Value2 = '  +4.04%  '
$ws.Cells.Item(29, 5).Value2 = '  -0.05%  '
$ws.Cells.Item(30, 5).Value2 = '  +1.31%  '
$ws.Cells.Item(31, 4).Value2 = '0.108'
$ws.Cells.Item(31, 5).Value2 = '  +0.18%  '
$ws.Cells.Item(32, 4).Value2 = '10.25'
$ws.Cells.Item(32, 5).Value2 = '  -2.59%  '
$ws.Cells.Item(33, 4).Value2 = '2.06'
$ws.Cells.Item(33, 5).Value2 = '  +0.25%  '
$ws.Cells.Item(34, 4).Value2 = '50.54'
$ws.Cells.Item(34, 5).Value2 = '  -0.44%  '
$ws.Cells.Item(35, 4).Value2 = '33.80'
$ws.Cells.Item(35, 5).Value2 = '  -0.39%  '
$ws.Cells.Item(36, 4).Value2 = '0.0448'
$ws.Cells.Item(37, 5).Value2 = '  -0.16%  '
$ws.Cells.Item(38, 4).Value2 = '3.30'
$ws.Cells.Item(38, 5).Value2 = '  +3.14%  '
$ws.Cells.Item(39, 4).Value2 = '0.290'
$ws.Cells.Item(39, 5).Value2 = '  +12.44%  '
$ws.Cells.Item(40, 4).Value2 = '16.88'
$ws.Cells.Item(40, 5).Value2 = '  +1.02%  '
$ws.Cells.Item(41, 5).Value2 = '  +1.28%  '
$ws.Cells.Item(42, 5).Value2 = '  -0.53%  '
$ws.Cells.Item(43, 2).Value2 = 'Monero'
$ws.Cells.Item(43, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(43, 4).Value2 = '126.64'
$ws.Cells.Item(43, 5).Value2 = '  +5.81%  '
$ws.Cells.Item(44, 2).Value2 = 'Stacks'
$ws.Cells.Item(44, 3).Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(44, 4).Value2 = '2.53'
$ws.Cells.Item(44, 5).Value2 = '  +1.89%  '
$ws.Cells.Item(45, 2).Value2 = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(45, 4).Value2 = '3.77'
$ws.Cells.Item(45, 5).Value2 = '  +5.90%  '
$ws.Cells.Item(46, 4).Value2 = '21.50'
$ws.Cells.Item(46, 5).Value2 = '  -0.04%  '
$ws.Cells.Item(47, 4).Value2 = '2.08'
$ws.Cells.Item(47, 5).Value2 = '  +2.63%  '
$ws.Cells.Item(48, 4).Value2 = '2.40'
$ws.Cells.Item(48, 5).Value2 = '  +2.50%  '
$ws.Cells.Item(49, 4).Value2 = '2.024.97'
$ws.Cells.Item(49, 5).Value2 = '  -0.64%  '
$ws.Cells.Item(50, 4).Value2 = '3.315.88'
$ws.Cells.Item(50, 5).Value2 = '  +2.07%  '
$ws.Cells.Item(51, 4).Value2 = '0.0316'
$ws.Cells.Item(51, 5).Value2 = '  -1.64%  '

# Restore original styles on the forced-text cells
$ws.Cells.Item(5, 4).Style = $origStyles["5_4"]
$ws.Cells.Item(6, 4).Style = $origStyles["6_4"]
$ws.Cells.Item(9, 4).Style = $origStyles["9_4"]
$ws.Cells.Item(10, 4).Style = $origStyles["10_4"]
$ws.Cells.Item(12, 4).Style = $origStyles["12_4"]
$ws.Cells.Item(14, 4).Style = $origStyles["14_4"]
$ws.Cells.Item(18, 4).Style = $origStyles["18_4"]
$ws.Cells.Item(20, 4).Style = $origStyles["20_4"]
$ws.Cells.Item(21, 4).Style = $origStyles["21_4"]
$ws.Cells.Item(23, 4).Style = $origStyles["23_4"]
$ws.Cells.Item(24, 4).Style = $origStyles["24_4"]
$ws.Cells.Item(26, 4).Style = $origStyles["26_4"]
$ws.Cells.Item(31, 4).Style = $origStyles["31_4"]
$ws.Cells.Item(32, 4).Style = $origStyles["32_4"]
$ws.Cells.Item(33, 4).Style = $origStyles["33_4"]
$ws.Cells.Item(34, 4).Style = $origStyles["34_4"]
$ws.Cells.Item(35, 4).Style = $origStyles["35_4"]
$ws.Cells.Item(36, 4).Style = $origStyles["36_4"]
$ws.Cells.Item(38, 4).Style = $origStyles["38_4"]
$ws.Cells.Item(39, 4).Style = $origStyles["39_4"]
$ws.Cells.Item(40, 4).Style = $origStyles["40_4"]
$ws.Cells.Item(43, 4).Style = $origStyles["43_4"]
$ws.Cells.Item(44, 4).Style = $origStyles["44_4"]
$ws.Cells.Item(45, 4).Style = $origStyles["45_4"]
$ws.Cells.Item(46, 4).Style = $origStyles["46_4"]
$ws.Cells.Item(47, 4).Style = $origStyles["47_4"]
$ws.Cells.Item(48, 4).Style = $origStyles["48_4"]
$ws.Cells.Item(51, 4).Style = $origStyles["51_4"]
